# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new value } for column F updates.
$updates = @{
    "展览" = @{
        3  = 527
        5  = 80
        6  = 31
        9  = 1121
        10 = 15584
        12 = 169
        14 = 6147
        17 = 64
        18 = 5
        19 = 110
        24 = 11
        27 = 861
        28 = 24
        31 = 11014
        34 = 112
    }
    "全部类型" = @{
        4  = 527
        6  = 80
        7  = 31
        10 = 1121
        11 = 15584
        13 = 169
        15 = 6147
        18 = 64
        19 = 5
        20 = 110
        25 = 11
        28 = 861
        29 = 24
        33 = 11014
        36 = 112
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
